$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2275
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2487.5
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7462.5
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -7798.5

$ws.Range("H112").Value = 1999.6666
$ws.Range("I112").Value = 2000
$ws.Range("K112").Value = 6000
$ws.Range("M112").Value = -4892

$ws.Range("H116").Value = 16286.667
$ws.Range("I116").Value = 16286.667
$ws.Range("K116").Value = 16286.667
$ws.Range("M116").Value = -12844.667

$ws.Range("H132").Value = 19336.818
$ws.Range("I132").Value = 21050.5
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 63151.5
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -60621.5
$ws.Range("N132").Value = -11660

$ws.Range("H138").Value = 2118.1052
$ws.Range("I138").Value = 1422.3636
$ws.Range("J138").Value = 3074.75
$ws.Range("K138").Value = 4267.0908
$ws.Range("L138").Value = 9224.25
$ws.Range("M138").Value = 872.9092000000001
$ws.Range("N138").Value = -19504.25

$ws.Range("H141").Value = 2592.7144
$ws.Range("I141").Value = 1691.5
$ws.Range("K141").Value = 5074.5
$ws.Range("M141").Value = 105.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 464.66666
$ws.Range("I25").Value = 197
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 197
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = 205
$ws.Range("N25").Value = -1804

$ws.Range("H74").Value = 2421.0667
$ws.Range("I74").Value = 2421.0667
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2421.0667
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1547.0667
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 2421.0667
$ws.Range("I77").Value = 2421.0667
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12105.3335
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7737.333499999999
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3898.926
$ws.Range("I107").Value = 1351.4667
$ws.Range("K107").Value = 1351.4667
$ws.Range("M107").Value = 568.5333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7488.154
$ws.Range("I31").Value = 4056.5
$ws.Range("K31").Value = 4056.5
$ws.Range("M31").Value = -3761.5

$ws.Range("H34").Value = 7488.154
$ws.Range("I34").Value = 4056.5
$ws.Range("K34").Value = 4056.5
$ws.Range("M34").Value = -3854.5

$ws.Range("H99").Value = 4918.857
$ws.Range("J99").Value = 5360.5
$ws.Range("L99").Value = 5360.5
$ws.Range("N99").Value = -8356.5

$ws.Range("H105").Value = 2022.5
$ws.Range("I105").Value = 2022.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2022.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -275.5
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 345.2143
$ws.Range("I107").Value = 403.5
$ws.Range("J107").Value = 199.5
$ws.Range("K107").Value = 403.5
$ws.Range("L107").Value = 199.5
$ws.Range("M107").Value = 1516.5
$ws.Range("N107").Value = -4039.5

$ws.Range("H126").Value = 4918.857
$ws.Range("J126").Value = 5360.5
$ws.Range("L126").Value = 16081.5
$ws.Range("N126").Value = -21021.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 289.57144
$ws.Range("I121").Value = 171.16667
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 513.50001
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 796.49999
$ws.Range("N121").Value = -5620

$ws.Range("H131").Value = 2052
$ws.Range("I131").Value = 2052
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6156
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1116
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2215.4614
$ws.Range("I80").Value = 1549.375
$ws.Range("J80").Value = 3281.2
$ws.Range("K80").Value = 1549.375
$ws.Range("L80").Value = 3281.2
$ws.Range("M80").Value = -551.375
$ws.Range("N80").Value = -5277.2

$ws.Range("H83").Value = 2215.4614
$ws.Range("I83").Value = 1549.375
$ws.Range("J83").Value = 3281.2
$ws.Range("K83").Value = 7746.875
$ws.Range("L83").Value = 16406
$ws.Range("M83").Value = -2754.875
$ws.Range("N83").Value = -26390

$ws.Range("H107").Value = 423.77777
$ws.Range("I107").Value = 419
$ws.Range("J107").Value = 433.33334
$ws.Range("K107").Value = 419
$ws.Range("L107").Value = 433.33334
$ws.Range("M107").Value = 1501
$ws.Range("N107").Value = -4273.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 17950
$ws.Range("I14").Value = 17950
$ws.Range("K14").Value = 17950
$ws.Range("M14").Value = -17778

$ws.Range("H22").Value = 1333.3334
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1333.3334
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1333.3334
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1923.3334

$ws.Range("H27").Value = 1333.3334
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1333.3334
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1333.3334
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1547.3334

$ws.Range("H136").Value = 2672.75
$ws.Range("I136").Value = 2672.75
$ws.Range("K136").Value = 8018.25
$ws.Range("M136").Value = -5468.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 7500
$ws.Range("J22").Value = 7500
$ws.Range("L22").Value = 7500
$ws.Range("N22").Value = -8086

$ws.Range("H31").Value = 1017
$ws.Range("I31").Value = 1017
$ws.Range("K31").Value = 1017
$ws.Range("M31").Value = -669

$ws.Range("H96").Value = 2199.3333
$ws.Range("I96").Value = 2149.25
$ws.Range("J96").Value = 2299.5
$ws.Range("K96").Value = 2149.25
$ws.Range("L96").Value = 2299.5
$ws.Range("M96").Value = -776.25
$ws.Range("N96").Value = -5045.5

$ws.Range("H132").Value = 3721.7646
$ws.Range("I132").Value = 2566.3333
$ws.Range("K132").Value = 7698.999899999999
$ws.Range("M132").Value = -5168.999899999999
